# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets
# to reflect the latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6595
$wsExhibit.Range("F4").Value = 105
$wsExhibit.Range("F5").Value = 148
$wsExhibit.Range("F7").Value = 79
$wsExhibit.Range("F8").Value = 587
$wsExhibit.Range("F9").Value = 45

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6595
$wsAll.Range("F5").Value = 105
$wsAll.Range("F6").Value = 148
$wsAll.Range("F9").Value = 79
$wsAll.Range("F10").Value = 587
$wsAll.Range("F11").Value = 45
